$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1318.3636
$ws.Range("I18").Value = 1183.6666
$ws.Range("J18").Value = 1480
$ws.Range("K18").Value = 1183.6666
$ws.Range("L18").Value = 1480
$ws.Range("M18").Value = -899.6666
$ws.Range("N18").Value = -2048
$ws.Range("H86").Value = 238948.75
$ws.Range("I86").Value = 1300
$ws.Range("J86").Value = 476597.5
$ws.Range("K86").Value = 1300
$ws.Range("L86").Value = 476597.5
$ws.Range("M86").Value = -177
$ws.Range("N86").Value = -478843.5
$ws.Range("H89").Value = 238948.75
$ws.Range("I89").Value = 1300
$ws.Range("J89").Value = 476597.5
$ws.Range("K89").Value = 6500
$ws.Range("L89").Value = 2382987.5
$ws.Range("M89").Value = -884
$ws.Range("N89").Value = -2394219.5
$ws.Range("H106").Value = 150000
$ws.Range("I106").Value = 150000
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 150000
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -149369
$ws.Range("H116").Value = 7099.7144
$ws.Range("I116").Value = 4999.5
$ws.Range("J116").Value = 7939.8
$ws.Range("K116").Value = 4999.5
$ws.Range("L116").Value = 7939.8
$ws.Range("M116").Value = -1557.5
$ws.Range("N116").Value = -14823.8
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 140.375
$ws.Range("I5").Value = 20.5
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 20.5
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 91.5
$ws.Range("H28").Value = 2437
$ws.Range("I28").Value = 2437
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2437
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2245
$ws.Range("H32").Value = 3503352
$ws.Range("I32").Value = 3503637.2
$ws.Range("J32").Value = 3500499.5
$ws.Range("K32").Value = 3503637.2
$ws.Range("L32").Value = 3500499.5
$ws.Range("M32").Value = -3503350.2
$ws.Range("N32").Value = -3501073.5
$ws.Range("H61").Value = 4000
$ws.Range("I61").Value = 4000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3788
$ws.Range("N61").Value = $null
$ws.Range("H80").Value = 80389.86
$ws.Range("I80").Value = 38999.5
$ws.Range("J80").Value = 96946
$ws.Range("K80").Value = 38999.5
$ws.Range("L80").Value = 96946
$ws.Range("M80").Value = -38001.5
$ws.Range("N80").Value = -98942
$ws.Range("H83").Value = 80389.86
$ws.Range("I83").Value = 38999.5
$ws.Range("J83").Value = 96946
$ws.Range("K83").Value = 116998.5
$ws.Range("L83").Value = 290838
$ws.Range("M83").Value = -112006.5
$ws.Range("N83").Value = -300822
$ws.Range("H99").Value = 2437
$ws.Range("I99").Value = 2437
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2437
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 558
$ws.Range("H122").Value = 3732.4
$ws.Range("I122").Value = 3331.5
$ws.Range("J122").Value = 3999.6667
$ws.Range("K122").Value = 9994.5
$ws.Range("L122").Value = 11999.0001
$ws.Range("M122").Value = -7544.5
$ws.Range("N122").Value = -16899.0001
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 140.375
$ws.Range("I4").Value = 20.5
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 20.5
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 94.5
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = $null
$ws.Range("M10").Value = 2000
$ws.Range("N10").Value = -2280
$ws.Range("H12").Value = 9999
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 9999
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = $null
$ws.Range("M12").Value = 9999
$ws.Range("N12").Value = -10335
$ws.Range("H20").Value = 1140.6
$ws.Range("I20").Value = 1140.6
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1140.6
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = -893.5999999999999
$ws.Range("H24").Value = 5000
$ws.Range("I24").Value = 5000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 5000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -4765
$ws.Range("H80").Value = 504.375
$ws.Range("I80").Value = 659.6
$ws.Range("J80").Value = 245.66667
$ws.Range("K80").Value = 659.6
$ws.Range("L80").Value = 245.66667
$ws.Range("M80").Value = 338.4
$ws.Range("N80").Value = -2241.66667
$ws.Range("H83").Value = 504.375
$ws.Range("I83").Value = 659.6
$ws.Range("J83").Value = 245.66667
$ws.Range("K83").Value = 3298
$ws.Range("L83").Value = 1228.33335
$ws.Range("M83").Value = 1694
$ws.Range("N83").Value = -11212.33335
$ws.Range("H86").Value = 1153.7142
$ws.Range("I86").Value = 997.2143
$ws.Range("J86").Value = 1466.7142
$ws.Range("K86").Value = 997.2143
$ws.Range("L86").Value = 1466.7142
$ws.Range("M86").Value = 125.7857
$ws.Range("H89").Value = 1153.7142
$ws.Range("I89").Value = 997.2143
$ws.Range("J89").Value = 1466.7142
$ws.Range("K89").Value = 4986.0715
$ws.Range("L89").Value = 7333.571
$ws.Range("M89").Value = 629.9285
$ws.Range("H94").Value = 1668.5238
$ws.Range("I94").Value = 1528.421
$ws.Range("J94").Value = 2999.5
$ws.Range("K94").Value = 1528.421
$ws.Range("L94").Value = 2999.5
$ws.Range("M94").Value = -1077.421
$ws.Range("H99").Value = 1847.2858
$ws.Range("I99").Value = 1995
$ws.Range("J99").Value = 1736.5
$ws.Range("K99").Value = 1995
$ws.Range("L99").Value = 1736.5
$ws.Range("M99").Value = -497
$ws.Range("N99").Value = -4732.5
$ws.Range("H132").Value = 74999.5
$ws.Range("I132").Value = 70000
$ws.Range("J132").Value = 79999
$ws.Range("K132").Value = 70000
$ws.Range("L132").Value = 79999
$ws.Range("M132").Value = -64940
$ws.Range("N132").Value = -90119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3398.6
$ws.Range("I16").Value = 4996.5
$ws.Range("J16").Value = 2333.3333
$ws.Range("K16").Value = 4996.5
$ws.Range("L16").Value = 2333.3333
$ws.Range("M16").Value = -4709.5
$ws.Range("N16").Value = -2907.3333
$ws.Range("H29").Value = 19990
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 19990
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = $null
$ws.Range("M29").Value = 19990
$ws.Range("N29").Value = -20576
$ws.Range("H86").Value = 10432
$ws.Range("I86").Value = 5650
$ws.Range("J86").Value = 19996
$ws.Range("K86").Value = 5650
$ws.Range("L86").Value = 19996
$ws.Range("M86").Value = -4527
$ws.Range("N86").Value = -22242
$ws.Range("H89").Value = 10432
$ws.Range("I89").Value = 5650
$ws.Range("J89").Value = 19996
$ws.Range("K89").Value = 28250
$ws.Range("L89").Value = 99980
$ws.Range("M89").Value = -22634
$ws.Range("N89").Value = -111212
$ws.Range("H94").Value = 114526.4
$ws.Range("I94").Value = 223564.2
$ws.Range("J94").Value = 5488.6
$ws.Range("K94").Value = 223564.2
$ws.Range("L94").Value = 5488.6
$ws.Range("M94").Value = -223113.2
$ws.Range("N94").Value = -6390.6
$ws.Range("H113").Value = 3398.6
$ws.Range("I113").Value = 4996.5
$ws.Range("J113").Value = 2333.3333
$ws.Range("K113").Value = 4996.5
$ws.Range("L113").Value = 2333.3333
$ws.Range("M113").Value = -2826.5
$ws.Range("N113").Value = -6673.3333
$ws.Range("H132").Value = 5945.357
$ws.Range("I132").Value = 6113.364
$ws.Range("J132").Value = 5329.3335
$ws.Range("K132").Value = 18340.092
$ws.Range("L132").Value = 15988.0005
$ws.Range("M132").Value = -15810.092
$ws.Range("N132").Value = -21048.0005
$ws.Range("H134").Value = 2460.8823
$ws.Range("I134").Value = 2427.1875
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 7281.5625
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4746.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 168.09091
$ws.Range("I12").Value = 235
$ws.Range("J12").Value = 51
$ws.Range("K12").Value = 705
$ws.Range("L12").Value = 153
$ws.Range("M12").Value = -532
$ws.Range("N12").Value = -499
$ws.Range("H33").Value = 510
$ws.Range("I33").Value = 220
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 1320
$ws.Range("L33").Value = 4800
$ws.Range("M33").Value = -1037
$ws.Range("N33").Value = -5366
$ws.Range("H86").Value = 150
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 150
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 450
$ws.Range("N86").Value = -2822
$ws.Range("H89").Value = 150
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 150
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 1350
$ws.Range("N89").Value = -13206

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 548.5
$ws.Range("I97").Value = 548.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 548.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -52.5
$ws.Range("H113").Value = 1524.5
$ws.Range("I113").Value = 287
$ws.Range("J113").Value = 3999.5
$ws.Range("K113").Value = 287
$ws.Range("L113").Value = 3999.5
$ws.Range("M113").Value = 1883
$ws.Range("H122").Value = 3033.3635

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2469.4666
$ws.Range("I46").Value = 2053.4285
$ws.Range("J46").Value = 2833.5
$ws.Range("K46").Value = 2053.4285
$ws.Range("L46").Value = 2833.5
$ws.Range("M46").Value = -1865.4285
$ws.Range("N46").Value = -3209.5
$ws.Range("H55").Value = 1525.5
$ws.Range("I55").Value = 1511.3334
$ws.Range("J55").Value = 1536.125
$ws.Range("K55").Value = 1511.3334
$ws.Range("L55").Value = 1536.125
$ws.Range("M55").Value = -1338.3334
$ws.Range("N55").Value = -1882.125
$ws.Range("H132").Value = 3013.1904
$ws.Range("I132").Value = 2865.4443
$ws.Range("J132").Value = 3899.6667
$ws.Range("K132").Value = 8596.332900000001
$ws.Range("L132").Value = 11699.0001
$ws.Range("M132").Value = -6066.332900000001
$ws.Range("H136").Value = 2589.6667
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2884.5
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 8653.5
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -13753.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6250624.5
$ws.Range("I100").Value = 7692873
$ws.Range("J100").Value = 881
$ws.Range("K100").Value = 15385746
$ws.Range("L100").Value = 1762
$ws.Range("M100").Value = -15385205
$ws.Range("N100").Value = -2844
